$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 108, pushing the existing
# rows 108-117 down to 110-119 (with their formatting/content intact).
$ws.Rows.Item(108).Insert()
$ws.Rows.Item(108).Insert()

# Fill in the two new rows (108 and 109) with the new data records.
$newRows = @(
    @{ Row = 108; D = 44610; K = "Red Globe";          M = 50; N = 10000; O = 10000; P = 10000; S = 556 },
    @{ Row = 109; D = 44610; K = "Superior Seedless";  M = 50; N = 10000; O = 10000; P = 10000; S = 556 }
)

foreach ($rec in $newRows) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value2  = 7
    $ws.Cells.Item($r, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
    $ws.Cells.Item($r, 3).Value2  = "Ñuble"
    $ws.Cells.Item($r, 4).Value2  = $rec.D
    $ws.Cells.Item($r, 5).Value2  = 16
    $ws.Cells.Item($r, 6).Value2  = "Fruta"
    $ws.Cells.Item($r, 7).Value2  = 100109
    $ws.Cells.Item($r, 8).Value2  = "Uva"
    $ws.Cells.Item($r, 9).Value2  = 100109001
    $ws.Cells.Item($r, 10).Value2 = "Uva"
    $ws.Cells.Item($r, 11).Value2 = $rec.K
    $ws.Cells.Item($r, 12).Value2 = "Primera"
    $ws.Cells.Item($r, 13).Value2 = $rec.M
    $ws.Cells.Item($r, 14).Value2 = $rec.N
    $ws.Cells.Item($r, 15).Value2 = $rec.O
    $ws.Cells.Item($r, 16).Value2 = $rec.P
    $ws.Cells.Item($r, 17).Value2 = "$/bandeja 18 kilos"
    $ws.Cells.Item($r, 18).Value2 = "Región de O'Higgins"
    $ws.Cells.Item($r, 19).Value2 = $rec.S
    $ws.Cells.Item($r, 20).Value2 = 18
}
